$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CART)
$ws.Range("B2").Value = 0.9457364341085271
$ws.Range("C2").Value = 0.9298245614035088
$ws.Range("D2").Value = 0.9177489177489178
$ws.Range("E2").Value = 0.9237472766884531

# Row 3 (C4.5)
$ws.Range("B3").Value = 0.9364341085271318
$ws.Range("C3").Value = 0.9166666666666666
$ws.Range("D3").Value = 0.9047619047619048
$ws.Range("E3").Value = 0.9106753812636166

# Row 8 (ExtraTrees)
$ws.Range("B8").Value = 0.8682170542635659
$ws.Range("C8").Value = 0.7280701754385965
$ws.Range("D8").Value = 0.8783068783068783
$ws.Range("E8").Value = 0.7961630695443646
